$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

# C1: add a border (top+bottom thin) matching the existing border def used for
# the "0" group header box under the merged B1:D1 cell.
$c1_C1 = $ws1.Range("C1")
$c1_C1.Style = "Normal"
$c1_C1.Borders.LineStyle = 1
$c1_C1.Borders.Item(7).LineStyle = -4142
$c1_C1.Borders.Item(10).LineStyle = -4142

# D1: add a border (top+bottom+right thin)
$c1_D1 = $ws1.Range("D1")
$c1_D1.Style = "Normal"
$c1_D1.Borders.LineStyle = 1
$c1_D1.Borders.Item(7).LineStyle = -4142

# C2: rename "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

# C1 / D1 (first group header box)
$c2_C1 = $ws2.Range("C1")
$c2_C1.Style = "Normal"
$c2_C1.Borders.LineStyle = 1
$c2_C1.Borders.Item(7).LineStyle = -4142
$c2_C1.Borders.Item(10).LineStyle = -4142

$c2_D1 = $ws2.Range("D1")
$c2_D1.Style = "Normal"
$c2_D1.Borders.LineStyle = 1
$c2_D1.Borders.Item(7).LineStyle = -4142

# F1 / G1 (second group header box)
$c2_F1 = $ws2.Range("F1")
$c2_F1.Style = "Normal"
$c2_F1.Borders.LineStyle = 1
$c2_F1.Borders.Item(7).LineStyle = -4142
$c2_F1.Borders.Item(10).LineStyle = -4142

$c2_G1 = $ws2.Range("G1")
$c2_G1.Style = "Normal"
$c2_G1.Borders.LineStyle = 1
$c2_G1.Borders.Item(7).LineStyle = -4142

# C2 / F2: rename "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5: remove the empty inline-string cell entirely
$ws2.Range("G5").ClearContents()
